$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 379.9565
$ws.Range("I33").Value = 381.6842
$ws.Range("J33").Value = 371.75
$ws.Range("K33").Value = 381.6842
$ws.Range("L33").Value = 371.75
$ws.Range("M33").Value = -152.6842
$ws.Range("N33").Value = -829.75
$ws.Range("H69").Value = 3602.889
$ws.Range("I69").Value = 3404.3333
$ws.Range("J69").Value = 4000
$ws.Range("K69").Value = 10212.9999
$ws.Range("L69").Value = 12000
$ws.Range("M69").Value = -9338.999899999999
$ws.Range("N69").Value = -13748
$ws.Range("H72").Value = 3602.889
$ws.Range("I72").Value = 3404.3333
$ws.Range("J72").Value = 4000
$ws.Range("K72").Value = 30638.9997
$ws.Range("L72").Value = 36000
$ws.Range("M72").Value = -26270.9997
$ws.Range("N72").Value = -44736
$ws.Range("H132").Value = 5438752.5
$ws.Range("I132").Value = 2530.9429
$ws.Range("J132").Value = 22735822
$ws.Range("K132").Value = 7592.8287
$ws.Range("L132").Value = 68207466
$ws.Range("M132").Value = -5062.8287
$ws.Range("N132").Value = -68212526
$ws.Range("H137").Value = 1095.875
$ws.Range("I137").Value = 1152.6471
$ws.Range("J137").Value = 958
$ws.Range("K137").Value = 3457.9413
$ws.Range("L137").Value = 2874
$ws.Range("M137").Value = -907.9412999999995
$ws.Range("N137").Value = -7974
$ws.Range("H138").Value = 4568708
$ws.Range("I138").Value = 9525489
$ws.Range("J138").Value = 3251.658
$ws.Range("K138").Value = 28576467
$ws.Range("L138").Value = 9754.974
$ws.Range("M138").Value = -28571327
$ws.Range("N138").Value = -20034.974
$ws.Range("H141").Value = 1532.6086
$ws.Range("I141").Value = 1330
$ws.Range("J141").Value = 2883.3333
$ws.Range("K141").Value = 3990
$ws.Range("L141").Value = 8649.999899999999
$ws.Range("M141").Value = 1190
$ws.Range("N141").Value = -19009.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1630
$ws.Range("I2").Value = 1556
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 1556
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -1443
$ws.Range("N2").Value = -2226
$ws.Range("H45").Value = 6496463.5
$ws.Range("I45").Value = 7578707.5
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 7578707.5
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -7578330.5
$ws.Range("N45").Value = -3754
$ws.Range("H74").Value = 1673.1923
$ws.Range("I74").Value = 1373.5
$ws.Range("K74").Value = 1373.5
$ws.Range("M74").Value = -499.5
$ws.Range("H77").Value = 1673.1923
$ws.Range("I77").Value = 1373.5
$ws.Range("K77").Value = 6867.5
$ws.Range("M77").Value = -2499.5
$ws.Range("H102").Value = 1014.5455
$ws.Range("I102").Value = 1028.8889
$ws.Range("K102").Value = 1028.8889
$ws.Range("M102").Value = 593.1111000000001
$ws.Range("H116").Value = 1630
$ws.Range("I116").Value = 1556
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 1556
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 738
$ws.Range("N116").Value = -6588
$ws.Range("H132").Value = 5103646
$ws.Range("I132").Value = 6411679
$ws.Range("J132").Value = 2317.1
$ws.Range("K132").Value = 19235037
$ws.Range("L132").Value = 6951.299999999999
$ws.Range("M132").Value = -19232507
$ws.Range("N132").Value = -12011.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1630
$ws.Range("I3").Value = 1556
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 1556
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -1442
$ws.Range("N3").Value = -2228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2151.2026
$ws.Range("I31").Value = 1957.7894
$ws.Range("J31").Value = 2355.361
$ws.Range("K31").Value = 1957.7894
$ws.Range("L31").Value = 2355.361
$ws.Range("M31").Value = -1662.7894
$ws.Range("N31").Value = -2945.361
$ws.Range("H34").Value = 2151.2026
$ws.Range("I34").Value = 1957.7894
$ws.Range("J34").Value = 2355.361
$ws.Range("K34").Value = 1957.7894
$ws.Range("L34").Value = 2355.361
$ws.Range("M34").Value = -1755.7894
$ws.Range("N34").Value = -2759.361

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 749.56525
$ws.Range("I113").Value = 626.5714
$ws.Range("J113").Value = 803.375
$ws.Range("K113").Value = 1879.7142
$ws.Range("L113").Value = 2410.125
$ws.Range("M113").Value = 290.2857999999999
$ws.Range("N113").Value = -6750.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12990.772
$ws.Range("I70").Value = 25822.111
$ws.Range("J70").Value = 4107.5386
$ws.Range("K70").Value = 25822.111
$ws.Range("L70").Value = 4107.5386
$ws.Range("M70").Value = -25552.111
$ws.Range("N70").Value = -4647.5386
$ws.Range("H73").Value = 12990.772
$ws.Range("I73").Value = 25822.111
$ws.Range("J73").Value = 4107.5386
$ws.Range("K73").Value = 25822.111
$ws.Range("L73").Value = 4107.5386
$ws.Range("M73").Value = -24886.111
$ws.Range("N73").Value = -5979.5386
$ws.Range("H80").Value = 18521386
$ws.Range("I80").Value = 41668990
$ws.Range("J80").Value = 3300
$ws.Range("K80").Value = 41668990
$ws.Range("L80").Value = 3300
$ws.Range("M80").Value = -41667992
$ws.Range("N80").Value = -5296
$ws.Range("H83").Value = 18521386
$ws.Range("I83").Value = 41668990
$ws.Range("J83").Value = 3300
$ws.Range("K83").Value = 208344950
$ws.Range("L83").Value = 16500
$ws.Range("M83").Value = -208339958
$ws.Range("N83").Value = -26484
$ws.Range("H88").Value = 38270
$ws.Range("J88").Value = 38270
$ws.Range("L88").Value = 38270
$ws.Range("N88").Value = -39172
$ws.Range("H91").Value = 38270
$ws.Range("J91").Value = 38270
$ws.Range("L91").Value = 38270
$ws.Range("N91").Value = -41390
$ws.Range("H97").Value = 1087.2916
$ws.Range("I97").Value = 1045.5294
$ws.Range("J97").Value = 1188.7142
$ws.Range("K97").Value = 1045.5294
$ws.Range("L97").Value = 1188.7142
$ws.Range("M97").Value = -549.5293999999999
$ws.Range("N97").Value = -2180.7142
$ws.Range("H102").Value = 3719.75
$ws.Range("I102").Value = 4435.3125
$ws.Range("J102").Value = 2288.625
$ws.Range("K102").Value = 4435.3125
$ws.Range("L102").Value = 2288.625
$ws.Range("M102").Value = -2813.3125
$ws.Range("N102").Value = -5532.625
$ws.Range("H104").Value = 50671
$ws.Range("J104").Value = 50671
$ws.Range("L104").Value = 50671
$ws.Range("N104").Value = -57659
$ws.Range("H132").Value = 5374.032
$ws.Range("I132").Value = 5741.731
$ws.Range("J132").Value = 3462
$ws.Range("K132").Value = 17225.193
$ws.Range("L132").Value = 10386
$ws.Range("M132").Value = -14695.193
$ws.Range("N132").Value = -15446

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 23700
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 23700
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 23700
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -24740
$ws.Range("H125").Value = 36337.5
$ws.Range("J125").Value = 36337.5
$ws.Range("L125").Value = 36337.5
$ws.Range("N125").Value = -46177.5
